$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert the new heading paragraph + a following blank paragraph in front
#    of the document's existing (sole) paragraph.
# ---------------------------------------------------------------------------
$enDash = [char]0x2013
$existing = $d.Paragraphs(1).Range
$existing.InsertBefore("Sentinel " + $enDash + " Level 400" + [char]13 + [char]13)

# Apply the "Heading 1" style to the newly created first paragraph. Word's
# COM host will materialise a skeleton "Heading1" paragraph style (and keep
# it linked to the built-in "Heading 1" Quick Style) the first time it is
# used in a document whose styles part does not yet define it.
$titlePara = $d.Paragraphs(1)
$titlePara.Style = "Heading 1"

# ---------------------------------------------------------------------------
# 2. Flesh out the auto-created "Heading1" paragraph style so that it mirrors
#    Word's standard Heading 1 definition (based on Normal, followed by
#    Normal, keepNext/keepLines, spacing before/after, accent1 colour at
#    32 half-points, linked character style, etc.)
# ---------------------------------------------------------------------------
$headingStyle = $d.Styles("Heading 1")
$headingStyle.NextParagraphStyle = "Normal"
$headingStyle.ParagraphFormat.SpaceBefore = 12
$headingStyle.ParagraphFormat.SpaceAfter = 0
$headingStyle.Font.Bold = $false
$headingStyle.Font.SizeBi = 16
$headingColor = $headingStyle.Font.TextColor
$headingColor.ObjectThemeColor = 4

# ---------------------------------------------------------------------------
# 3. Create + link the paired "Heading 1 Char" character style, matching the
#    built-in style gallery's definition.
# ---------------------------------------------------------------------------
$charStyle = $d.Styles.Add("Heading1Char", 2)
$headingStyle.LinkStyle = $charStyle
$charStyle.NameLocal = "Heading 1 Char"
$charStyle.BaseStyle = "DefaultParagraphFont"
$charStyle.LinkStyle = "Heading1"
$charStyle.Priority = 9
$charStyle.Font.Size = 16
$charStyle.Font.SizeBi = 16
$charColor = $charStyle.Font.TextColor
$charColor.ObjectThemeColor = 4

# Normalise the display name to match Word's own "heading 1" (lower-case)
# after every by-name lookup above is already done.
$headingStyle.NameLocal = "heading 1"
